# Updates the bilibili-show listing sheets ("展览" / "全部类型") to the
# newly scraped data: refresh a handful of "want-to-go" counters and
# splice in a newly-discovered event ("合肥·星光国潮动漫游戏嘉年华")
# right before the last (existing) row, pushing that last row down by one.

function Update-EventSheet($ws, $lastRow, $counterUpdates) {
    # --- 1. refresh "want to go" counters (column F) on existing rows ----
    foreach ($pair in $counterUpdates) {
        $ws.Cells.Item($pair[0], 6).Value = $pair[1]
    }

    # --- 2. make room for the newly discovered event right before the ---
    # --- last row, shifting the old last row one row down ---------------
    $ws.Rows.Item($lastRow).Insert()

    # column A carries direct formatting (bold / border / centred) -
    # copy it down from the row above so the new row matches the rest
    # of the table instead of inheriting the insert's blank format.
    $ws.Cells.Item($lastRow - 1, 1).Copy($ws.Cells.Item($lastRow, 1))

    # --- 3. populate the new row with the newly scraped event -----------
    $ws.Cells.Item($lastRow, 1).Value = $lastRow - 1

    # Column B holds plain date-text ("YYYY-MM-DD"); force text format
    # first so Excel doesn't silently reinterpret it as a date serial,
    # then drop back to the default style so no stray formatting sticks.
    $ws.Cells.Item($lastRow, 2).NumberFormat = "@"
    $ws.Cells.Item($lastRow, 2).Value = "2024-12-08"
    $ws.Cells.Item($lastRow, 2).Style = "Normal"

    $ws.Cells.Item($lastRow, 3).Value = "合肥·星光国潮动漫游戏嘉年华"
    $ws.Cells.Item($lastRow, 4).Value = "北二环与新蚌埠路交汇处 蓝金湾大酒店"
    $ws.Cells.Item($lastRow, 5).Value = "2024.12.08 10:00-12.08 17:00"
    $ws.Cells.Item($lastRow, 6).Value = 0
    $ws.Cells.Item($lastRow, 7).Value = 39.9
    $ws.Cells.Item($lastRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=93801"
    $ws.Cells.Item($lastRow, 9).Value = "//i0.hdslb.com/bfs/openplatform/202410/ubX6VZ841729253636894.png"

    # --- 4. the shifted (was-last) row keeps its own data, only its -----
    # --- sequence number and "want to go" counter move on ---------------
    $ws.Cells.Item($lastRow + 1, 1).Value = $lastRow
    $ws.Cells.Item($lastRow + 1, 6).Value = 59
}

$wb = $excel.ActiveWorkbook

# "展览" (exhibitions) — rows 1..10, new row spliced in before row 10
$wsExhibitions = $wb.Worksheets.Item(1)
$exhibitionCounters = @(
    , @(2, 452)
    , @(4, 68)
    , @(5, 5118)
    , @(7, 39)
    , @(8, 94)
    , @(9, 324)
)
Update-EventSheet $wsExhibitions 10 $exhibitionCounters

# "全部类型" (all types) — rows 1..15, new row spliced in before row 15
$wsAll = $wb.Worksheets.Item(4)
$allTypesCounters = @(
    , @(2, 452)
    , @(8, 68)
    , @(9, 5118)
    , @(11, 39)
    , @(12, 94)
    , @(14, 324)
)
Update-EventSheet $wsAll 15 $allTypesCounters
